$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E).
# Some Price values look like plain numbers (e.g. "0.517", "2.50"); those
# are written with a leading apostrophe so Excel keeps them as text
# instead of auto-converting them to numeric values, matching the
# original inline-string cell contents.
$ws.Range("D2").Value = '45.491.24'
$ws.Range("E2").Value = '  +3.32%  '
$ws.Range("D3").Value = '2.429.19'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''318.38'
$ws.Range("E5").Value = '  +3.46%  '
$ws.Range("D6").Value = '''102.93'
$ws.Range("E6").Value = '  +5.26%  '
$ws.Range("D7").Value = '''0.517'
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.531'
$ws.Range("E9").Value = '  +6.91%  '
$ws.Range("D10").Value = '''35.71'
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").Value = '''0.0804'
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D13").Value = '''18.15'
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").Value = '''7.07'
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").Value = '2.808.78'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '2.398.83'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '''0.844'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '45.353.56'
$ws.Range("E18").Value = '  +3.07%  '
$ws.Range("D19").Value = '''12.27'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").Value = '0.0₃0922'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").Value = '''68.88'
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").Value = '''244.58'
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '''2.50'
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '''25.58'
$ws.Range("E27").Value = '  +1.77%  '
$ws.Range("D28").Value = '''2.18'
$ws.Range("E28").Value = '  -0.96%  '
$ws.Range("D29").Value = '''9.59'
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").Value = '''49.26'
$ws.Range("E30").Value = '  +2.53%  '
$ws.Range("D31").Value = '''32.92'
$ws.Range("E31").Value = '  +1.11%  '
$ws.Range("D32").Value = '''20.31'
$ws.Range("E32").Value = '  +9.48%  '
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("D34").Value = '''5.22'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").Value = '''0.0769'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("D37").Value = '''1.88'
$ws.Range("E37").Value = '  -3.64%  '
$ws.Range("D38").Value = '''4.46'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("D40").Value = '''125.32'
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("E41").Value = '  -2.48%  '
$ws.Range("D42").Value = '''0.110'
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").Value = '''20.38'
$ws.Range("E43").Value = '  -3.78%  '
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").Value = '1.929.41'
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("E46").Value = '  -2.84%  '
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("E48").Value = '  +13.69%  '
$ws.Range("D49").Value = '''9.12'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").Value = '''76.85'
$ws.Range("E50").Value = '  +5.47%  '
$ws.Range("D51").Value = '''53.86'
$ws.Range("E51").Value = '  +1.64%  '
